# Add new "Meaning" rows (grade-span code lookups) to the DataMeaningType sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New data rows 295-300 (columns A, B, C) ---
$data = @(
    @("A", "School has elementary, middle, and high school grades"),
    @("E", "School has elementary schools grades (PK-5)"),
    @("H", "School has high schools grades (9-13)"),
    @("I", "Schools with elementary, middle and high school grades"),
    @("M", "School has middle schools grades (6-8)"),
    @("T", "School has middle and high school grades")
)

$startRow = 295
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = "Object"
}

# --- 2. Apply the small Arial 7pt font used for this block of rows.
#        Build the style once on a scratch cell (within the already-used
#        column range so the sheet dimension is not disturbed), stamp the
#        whole B295:B300 range with one PasteSpecial(Formats), then remove
#        the scratch cell completely. ---
$helperCell = $ws.Range("A1000")
$helperCell.Value = "fmt"
$hf = $helperCell.Font
$hf.Name = "Arial"
$hf.Size = 7
$hf.Color = 0

$helperCell.Copy()
$ws.Range("B295:B300").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$helperCell.Clear()

# --- 3. Cosmetic view-state updates to mirror the edited workbook ---
$ws.Range("A295:A300").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 287
$win.ScrollColumn = 1

# --- 4. Page setup: orientation portrait (as introduced in the diff) ---
$ps = $ws.PageSetup
$ps.Orientation = 1
